# Update "想去人数" (column F) counts on the "展览" and "全部类型" sheets
# to reflect the latest generated output.

$wb = $excel.ActiveWorkbook

# Map of sheet name -> list of (cell, newValue) updates
$updates = @{
    "展览"     = @{ "F2" = 1175; "F3" = 420; "F7" = 12300; "F11" = 148; "F12" = 12097; "F14" = 4671; "F16" = 56; "F22" = 166; "F23" = 71 }
    "全部类型" = @{ "F2" = 1175; "F3" = 420; "F9" = 12300; "F13" = 148; "F14" = 12097; "F16" = 4671; "F18" = 56; "F24" = 166; "F25" = 71 }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $cellUpdates = $updates[$sheetName]
    foreach ($cellRef in $cellUpdates.Keys) {
        $ws.Range($cellRef).Value = $cellUpdates[$cellRef]
    }
}
